# Fruta / hortaliza, semanal
# Weekly update: a new price observation (for 44474) is inserted ahead of the
# existing rows 17-19, which are pushed down to become rows 19-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 17. This shifts the current rows
# 17 ("Packham's Triumph"), 18 ("Winter Nelis") and 19 ("Packham's Triumph",
# bandeja) down to rows 19, 20 and 21 respectively, preserving their values
# and formatting (including the date-style on column D) automatically.
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()

# New row 17: updated "Packham's Triumph" observation.
$ws.Cells.Item(17, 1).Value = 1
$ws.Cells.Item(17, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(17, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(17, 4).Value = 44474
$ws.Cells.Item(17, 5).Value = 15
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100104
$ws.Cells.Item(17, 8).Value = "Frutos de pepita"
$ws.Cells.Item(17, 9).Value = 100104005
$ws.Cells.Item(17, 10).Value = "Pera"
$ws.Cells.Item(17, 11).Value = "Packham's Triumph"
$ws.Cells.Item(17, 12).Value = "Segunda"
$ws.Cells.Item(17, 13).Value = 270
$ws.Cells.Item(17, 14).Value = 18000
$ws.Cells.Item(17, 15).Value = 19000
$ws.Cells.Item(17, 16).Value = 18500
$ws.Cells.Item(17, 17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(17, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(17, 19).Value = 1028
$ws.Cells.Item(17, 20).Value = 18

# New row 18: updated "Winter Nelis" observation.
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(18, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(18, 4).Value = 44474
$ws.Cells.Item(18, 5).Value = 15
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100104
$ws.Cells.Item(18, 8).Value = "Frutos de pepita"
$ws.Cells.Item(18, 9).Value = 100104005
$ws.Cells.Item(18, 10).Value = "Pera"
$ws.Cells.Item(18, 11).Value = "Winter Nelis"
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 250
$ws.Cells.Item(18, 14).Value = 17000
$ws.Cells.Item(18, 15).Value = 18000
$ws.Cells.Item(18, 16).Value = 17500
$ws.Cells.Item(18, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(18, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(18, 19).Value = 972
$ws.Cells.Item(18, 20).Value = 18
